$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update measured leg values (this cascades automatically into the
# dependent formula cells H, I, J, K, L, M through recalculation)

# Row 4 (L1)
$ws.Range("C4").Value = 1430
$ws.Range("F4").Value = 1950

# Row 5 (L2)
$ws.Range("C5").Value = 1580
$ws.Range("D5").Value = 1570
$ws.Range("G5").Value = 2000

# Row 6 (L3)
$ws.Range("C6").Value = 1400
$ws.Range("D6").Value = 1680
$ws.Range("G6").Value = 2150

# Row 8 (R2)
$ws.Range("D8").Value = 1200
$ws.Range("G8").Value = 730

# Row 9 (R3)
$ws.Range("C9").Value = 1430
$ws.Range("D9").Value = 1230
$ws.Range("E9").Value = 1830
$ws.Range("F9").Value = 900

# Recalculate so dependent formula cells are refreshed
$excel.CalculateFullRebuild()

# Update the active selection on the sheet to match the saved view
$ws.Range("G11").Select()
